$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("sample_name") contains identifiable sample reference codes.
# These are being anonymized/replaced with generic "Sample N" labels.
$ws.Range("F2").Value = "Sample 1"
$ws.Range("F3").Value = "Sample 11"
$ws.Range("F4").Value = "Sample 14"
$ws.Range("F5").Value = "Sample 16"
$ws.Range("F6").Value = "Sample 22"
$ws.Range("F7").Value = "Sample 26"
$ws.Range("F8").Value = "Sample 47"
$ws.Range("F9").Value = "Sample 48"
$ws.Range("F10").Value = "Sample 88"
$ws.Range("F11").Value = "Sample 100"
